$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All 150 new strings in the order they appear in the shared string table:
# indices 0..147 -> I3..I150 (path then file list)
# index 148 -> I1 ("d2flats")
# index 149 -> I2 (long comment)
$values = @(
    'C:\Users\TwoIsOne\Desktop\5.22.19\WORKSPACE8\BAK\D2Flats\',
    'BLOOD1.png',
    'BLOOD2.png',
    'BLOOD3.png',
    'CEIL1_1.png',
    'CEIL1_2.png',
    'CEIL1_3.png',
    'CEIL3_1.png',
    'CEIL3_2.png',
    'CEIL3_3.png',
    'CEIL3_4.png',
    'CEIL3_5.png',
    'CEIL3_6.png',
    'CEIL4_1.png',
    'CEIL4_2.png',
    'CEIL4_3.png',
    'CEIL5_1.png',
    'CEIL5_2.png',
    'COMP01.png',
    'CONS1_1.png',
    'CONS1_5.png',
    'CONS1_7.png',
    'CRATOP1.png',
    'CRATOP2.png',
    'DEM1_1.png',
    'DEM1_2.png',
    'DEM1_3.png',
    'DEM1_4.png',
    'DEM1_5.png',
    'DEM1_6.png',
    'F_SKY1.png',
    'FLAT1.png',
    'FLAT1_1.png',
    'FLAT1_2.png',
    'FLAT1_3.png',
    'FLAT2.png',
    'FLAT3.png',
    'FLAT4.png',
    'FLAT5.png',
    'FLAT5_1.png',
    'FLAT5_2.png',
    'FLAT5_3.png',
    'FLAT5_4.png',
    'FLAT5_5.png',
    'FLAT5_6.png',
    'FLAT5_7.png',
    'FLAT5_8.png',
    'FLAT8.png',
    'FLAT9.png',
    'FLAT10.png',
    'FLAT14.png',
    'FLAT17.png',
    'FLAT18.png',
    'FLAT19.png',
    'FLAT20.png',
    'FLAT22.png',
    'FLAT23.png',
    'FLOOR0_1.png',
    'FLOOR0_2.png',
    'FLOOR0_3.png',
    'FLOOR0_5.png',
    'FLOOR0_6.png',
    'FLOOR0_7.png',
    'FLOOR1_1.png',
    'FLOOR1_6.png',
    'FLOOR1_7.png',
    'FLOOR3_3.png',
    'FLOOR4_1.png',
    'FLOOR4_5.png',
    'FLOOR4_6.png',
    'FLOOR4_8.png',
    'FLOOR5_1.png',
    'FLOOR5_2.png',
    'FLOOR5_3.png',
    'FLOOR5_4.png',
    'FLOOR6_1.png',
    'FLOOR6_2.png',
    'FLOOR7_1.png',
    'FLOOR7_2.png',
    'FWATER1.png',
    'FWATER2.png',
    'FWATER3.png',
    'FWATER4.png',
    'GATE1.png',
    'GATE2.png',
    'GATE3.png',
    'GATE4.png',
    'GRASS1.png',
    'GRASS2.png',
    'GRNLITE1.png',
    'GRNROCK.png',
    'LAVA1.png',
    'LAVA2.png',
    'LAVA3.png',
    'LAVA4.png',
    'MFLR8_1.png',
    'MFLR8_2.png',
    'MFLR8_3.png',
    'MFLR8_4.png',
    'NUKAGE1.png',
    'NUKAGE2.png',
    'NUKAGE3.png',
    'RROCK01.png',
    'RROCK02.png',
    'RROCK03.png',
    'RROCK04.png',
    'RROCK05.png',
    'RROCK06.png',
    'RROCK07.png',
    'RROCK08.png',
    'RROCK09.png',
    'RROCK10.png',
    'RROCK11.png',
    'RROCK12.png',
    'RROCK13.png',
    'RROCK14.png',
    'RROCK15.png',
    'RROCK16.png',
    'RROCK17.png',
    'RROCK18.png',
    'RROCK19.png',
    'RROCK20.png',
    'SFLR6_1.png',
    'SFLR6_4.png',
    'SFLR7_1.png',
    'SFLR7_4.png',
    'SLIME01.png',
    'SLIME02.png',
    'SLIME03.png',
    'SLIME04.png',
    'SLIME05.png',
    'SLIME06.png',
    'SLIME07.png',
    'SLIME08.png',
    'SLIME09.png',
    'SLIME10.png',
    'SLIME11.png',
    'SLIME12.png',
    'SLIME13.png',
    'SLIME14.png',
    'SLIME15.png',
    'SLIME16.png',
    'STEP1.png',
    'STEP2.png',
    'TLITE6_1.png',
    'TLITE6_4.png',
    'TLITE6_5.png',
    'TLITE6_6.png',
    'd2flats',
    'some excel magic could help to see what innate flats has been disamiguated, they will be much more common to need replacement in the wads'
)

for ($i = 0; $i -lt 148; $i++) {
    $row = $i + 3
    $ws.Cells.Item($row, 9).Value = $values[$i]
}

$ws.Range("I1").Value = $values[148]
$ws.Range("I2").Value = $values[149]

$ws.Columns.Item(9).ColumnWidth = 31.6

$ws.Range("I2").Select() | Out-Null
